$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 454, shifting rows 454-528 down to 455-529
$ws.Rows("454:454").Insert()

# Populate the new row 454 with the new data
$ws.Cells.Item(454, 1).Value = 4
$ws.Cells.Item(454, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(454, 3).Value = "Los Lagos"
$ws.Cells.Item(454, 4).Value = 45180
$ws.Cells.Item(454, 5).Value = 10
$ws.Cells.Item(454, 6).Value = 100112040
$ws.Cells.Item(454, 7).Value = "Cilantro"
$ws.Cells.Item(454, 8).Value = "Sin especificar"
$ws.Cells.Item(454, 9).Value = "Primera"
$ws.Cells.Item(454, 10).Value = 80
$ws.Cells.Item(454, 11).Value = 14000
$ws.Cells.Item(454, 12).Value = 14000
$ws.Cells.Item(454, 13).Value = 14000
$ws.Cells.Item(454, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(454, 15).Value = "Región Metropolitana"
$ws.Cells.Item(454, 16).Value = 389
$ws.Cells.Item(454, 17).Value = 36
$ws.Cells.Item(454, 18).Value = "Hortaliza"
